# Update the "想去人数" (column F) figures on the "展览" and "全部类型"
# sheets to match the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 3119
$ws1.Range("F3").Value  = 520
$ws1.Range("F4").Value  = 952
$ws1.Range("F5").Value  = 77
$ws1.Range("F6").Value  = 23
$ws1.Range("F9").Value  = 1102
$ws1.Range("F10").Value = 15417
$ws1.Range("F11").Value = 220
$ws1.Range("F12").Value = 163
$ws1.Range("F14").Value = 6080
$ws1.Range("F15").Value = 618
$ws1.Range("F16").Value = 101
$ws1.Range("F20").Value = 1260
$ws1.Range("F22").Value = 110
$ws1.Range("F27").Value = 850
$ws1.Range("F29").Value = 4982
$ws1.Range("F30").Value = 135
$ws1.Range("F31").Value = 10949
$ws1.Range("F32").Value = 1224
$ws1.Range("F35").Value = 153
$ws1.Range("F36").Value = 3784
$ws1.Range("F38").Value = 72

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 3119
$ws4.Range("F4").Value  = 520
$ws4.Range("F5").Value  = 953
$ws4.Range("F6").Value  = 77
$ws4.Range("F7").Value  = 23
$ws4.Range("F10").Value = 1102
$ws4.Range("F11").Value = 15417
$ws4.Range("F12").Value = 220
$ws4.Range("F13").Value = 163
$ws4.Range("F15").Value = 6080
$ws4.Range("F16").Value = 618
$ws4.Range("F17").Value = 101
$ws4.Range("F21").Value = 1260
$ws4.Range("F23").Value = 110
$ws4.Range("F28").Value = 850
$ws4.Range("F30").Value = 4982
$ws4.Range("F31").Value = 135
$ws4.Range("F33").Value = 10949
$ws4.Range("F34").Value = 1224
$ws4.Range("F37").Value = 153
$ws4.Range("F38").Value = 3784
$ws4.Range("F40").Value = 72

$wb.Save()
